$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Tfpi"
$ws.Cells.Item(2, 3).Value = "Vldlr"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 10.91244366666667
$ws.Cells.Item(2, 8).Value = 32.737331
$ws.Cells.Item(2, 9).Value = 0.2863847678890978
$ws.Cells.Item(2, 10).Value = 0.2863847678890977
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.3081963333333333
$ws.Cells.Item(2, 14).Value = 0.924589
$ws.Cells.Item(2, 15).Value = 0.09210955608663024
$ws.Cells.Item(2, 16).Value = 0.09210955608663024
$ws.Cells.Item(2, 17).Value = 3.363175125773223
$ws.Cells.Item(2, 18).Value = 30.268576131959
$ws.Cells.Item(2, 19).Value = 0.02637877384023744
$ws.Cells.Item(2, 20).Value = 0.02637877384023743

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Tfpi"
$ws.Cells.Item(3, 3).Value = "Vldlr"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 10.91244366666667
$ws.Cells.Item(3, 8).Value = 32.737331
$ws.Cells.Item(3, 9).Value = 0.2863847678890978
$ws.Cells.Item(3, 10).Value = 0.2863847678890977
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.207039333333333
$ws.Cells.Item(3, 14).Value = 6.621118
$ws.Cells.Item(3, 15).Value = 0.65961009678592
$ws.Cells.Item(3, 16).Value = 0.6596100967859201
$ws.Cells.Item(3, 17).Value = 24.08419239511755
$ws.Cells.Item(3, 18).Value = 216.757731556058
$ws.Cells.Item(3, 19).Value = 0.188902284465341
$ws.Cells.Item(3, 20).Value = 0.188902284465341

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Tfpi"
$ws.Cells.Item(4, 3).Value = "Vldlr"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 10.91244366666667
$ws.Cells.Item(4, 8).Value = 32.737331
$ws.Cells.Item(4, 9).Value = 0.2863847678890978
$ws.Cells.Item(4, 10).Value = 0.2863847678890977
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.83074
$ws.Cells.Item(4, 14).Value = 2.49222
$ws.Cells.Item(4, 15).Value = 0.2482803471274497
$ws.Cells.Item(4, 16).Value = 0.2482803471274497
$ws.Cells.Item(4, 17).Value = 9.065403451646667
$ws.Cells.Item(4, 18).Value = 81.58863106481999
$ws.Cells.Item(4, 19).Value = 0.0711037095835193
$ws.Cells.Item(4, 20).Value = 0.0711037095835193

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Tfpi"
$ws.Cells.Item(5, 3).Value = "Vldlr"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 17.16042866666666
$ws.Cells.Item(5, 8).Value = 51.481286
$ws.Cells.Item(5, 9).Value = 0.4503560825328814
$ws.Cells.Item(5, 10).Value = 0.4503560825328813
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.3081963333333333
$ws.Cells.Item(5, 14).Value = 0.924589
$ws.Cells.Item(5, 15).Value = 0.09210955608663024
$ws.Cells.Item(5, 16).Value = 0.09210955608663024
$ws.Cells.Item(5, 17).Value = 5.288781193494889
$ws.Cells.Item(5, 18).Value = 47.59903074145399
$ws.Cells.Item(5, 19).Value = 0.04148209884301752
$ws.Cells.Item(5, 20).Value = 0.04148209884301751

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Tfpi"
$ws.Cells.Item(6, 3).Value = "Vldlr"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 17.16042866666666
$ws.Cells.Item(6, 8).Value = 51.481286
$ws.Cells.Item(6, 9).Value = 0.4503560825328814
$ws.Cells.Item(6, 10).Value = 0.4503560825328813
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.207039333333333
$ws.Cells.Item(6, 14).Value = 6.621118
$ws.Cells.Item(6, 15).Value = 0.65961009678592
$ws.Cells.Item(6, 16).Value = 0.6596100967859201
$ws.Cells.Item(6, 17).Value = 37.87374104419422
$ws.Cells.Item(6, 18).Value = 340.863669397748
$ws.Cells.Item(6, 19).Value = 0.2970594191876417
$ws.Cells.Item(6, 20).Value = 0.2970594191876417

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Tfpi"
$ws.Cells.Item(7, 3).Value = "Vldlr"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 17.16042866666666
$ws.Cells.Item(7, 8).Value = 51.481286
$ws.Cells.Item(7, 9).Value = 0.4503560825328814
$ws.Cells.Item(7, 10).Value = 0.4503560825328813
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.83074
$ws.Cells.Item(7, 14).Value = 2.49222
$ws.Cells.Item(7, 15).Value = 0.2482803471274497
$ws.Cells.Item(7, 16).Value = 0.2482803471274497
$ws.Cells.Item(7, 17).Value = 14.25585451054667
$ws.Cells.Item(7, 18).Value = 128.30269059492
$ws.Cells.Item(7, 19).Value = 0.1118145645022222
$ws.Cells.Item(7, 20).Value = 0.1118145645022222

# Row 8
$ws.Cells.Item(8, 1).Value = "M1"
$ws.Cells.Item(8, 2).Value = "Tfpi"
$ws.Cells.Item(8, 3).Value = "Vldlr"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 4.33904
$ws.Cells.Item(8, 8).Value = 13.01712
$ws.Cells.Item(8, 9).Value = 0.1138732076168498
$ws.Cells.Item(8, 10).Value = 0.1138732076168498
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.3081963333333333
$ws.Cells.Item(8, 14).Value = 0.924589
$ws.Cells.Item(8, 15).Value = 0.09210955608663024
$ws.Cells.Item(8, 16).Value = 0.09210955608663024
$ws.Cells.Item(8, 17).Value = 1.337276218186667
$ws.Cells.Item(8, 18).Value = 12.03548596368
$ws.Cells.Item(8, 19).Value = 0.01048881060374872
$ws.Cells.Item(8, 20).Value = 0.01048881060374871

# Row 9
$ws.Cells.Item(9, 1).Value = "M1"
$ws.Cells.Item(9, 2).Value = "Tfpi"
$ws.Cells.Item(9, 3).Value = "Vldlr"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 4.33904
$ws.Cells.Item(9, 8).Value = 13.01712
$ws.Cells.Item(9, 9).Value = 0.1138732076168498
$ws.Cells.Item(9, 10).Value = 0.1138732076168498
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.207039333333333
$ws.Cells.Item(9, 14).Value = 6.621118
$ws.Cells.Item(9, 15).Value = 0.65961009678592
$ws.Cells.Item(9, 16).Value = 0.6596100967859201
$ws.Cells.Item(9, 17).Value = 9.576431948906666
$ws.Cells.Item(9, 18).Value = 86.18788754016001
$ws.Cells.Item(9, 19).Value = 0.07511191749747345
$ws.Cells.Item(9, 20).Value = 0.07511191749747345

# Row 10
$ws.Cells.Item(10, 1).Value = "M1"
$ws.Cells.Item(10, 2).Value = "Tfpi"
$ws.Cells.Item(10, 3).Value = "Vldlr"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 4.33904
$ws.Cells.Item(10, 8).Value = 13.01712
$ws.Cells.Item(10, 9).Value = 0.1138732076168498
$ws.Cells.Item(10, 10).Value = 0.1138732076168498
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.83074
$ws.Cells.Item(10, 14).Value = 2.49222
$ws.Cells.Item(10, 15).Value = 0.2482803471274497
$ws.Cells.Item(10, 16).Value = 0.2482803471274497
$ws.Cells.Item(10, 17).Value = 3.6046140896
$ws.Cells.Item(10, 18).Value = 32.4415268064
$ws.Cells.Item(10, 19).Value = 0.02827247951562761
$ws.Cells.Item(10, 20).Value = 0.02827247951562761

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Tfpi"
$ws.Cells.Item(11, 3).Value = "Vldlr"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 3.892394333333333
$ws.Cells.Item(11, 8).Value = 11.677183
$ws.Cells.Item(11, 9).Value = 0.1021514961941619
$ws.Cells.Item(11, 10).Value = 0.1021514961941619
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.3081963333333333
$ws.Cells.Item(11, 14).Value = 0.924589
$ws.Cells.Item(11, 15).Value = 0.09210955608663024
$ws.Cells.Item(11, 16).Value = 0.09210955608663024
$ws.Cells.Item(11, 17).Value = 1.199621661420778
$ws.Cells.Item(11, 18).Value = 10.796594952787
$ws.Cells.Item(11, 19).Value = 0.009409128968029352
$ws.Cells.Item(11, 20).Value = 0.009409128968029349

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Tfpi"
$ws.Cells.Item(12, 3).Value = "Vldlr"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 3.892394333333333
$ws.Cells.Item(12, 8).Value = 11.677183
$ws.Cells.Item(12, 9).Value = 0.1021514961941619
$ws.Cells.Item(12, 10).Value = 0.1021514961941619
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 2.207039333333333
$ws.Cells.Item(12, 14).Value = 6.621118
$ws.Cells.Item(12, 15).Value = 0.65961009678592
$ws.Cells.Item(12, 16).Value = 0.6596100967859201
$ws.Cells.Item(12, 17).Value = 8.590667394510444
$ws.Cells.Item(12, 18).Value = 77.316006550594
$ws.Cells.Item(12, 19).Value = 0.06738015829145767
$ws.Cells.Item(12, 20).Value = 0.06738015829145766

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Tfpi"
$ws.Cells.Item(13, 3).Value = "Vldlr"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 3.892394333333333
$ws.Cells.Item(13, 8).Value = 11.677183
$ws.Cells.Item(13, 9).Value = 0.1021514961941619
$ws.Cells.Item(13, 10).Value = 0.1021514961941619
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.83074
$ws.Cells.Item(13, 14).Value = 2.49222
$ws.Cells.Item(13, 15).Value = 0.2482803471274497
$ws.Cells.Item(13, 16).Value = 0.2482803471274497
$ws.Cells.Item(13, 17).Value = 3.233567668473333
$ws.Cells.Item(13, 18).Value = 29.10210901626
$ws.Cells.Item(13, 19).Value = 0.02536220893467487
$ws.Cells.Item(13, 20).Value = 0.02536220893467487

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Tfpi"
$ws.Cells.Item(14, 3).Value = "Vldlr"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 1.799827666666667
$ws.Cells.Item(14, 8).Value = 5.399483
$ws.Cells.Item(14, 9).Value = 0.04723444576700921
$ws.Cells.Item(14, 10).Value = 0.0472344457670092
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 0.3081963333333333
$ws.Cells.Item(14, 14).Value = 0.924589
$ws.Cells.Item(14, 15).Value = 0.09210955608663024
$ws.Cells.Item(14, 16).Value = 0.09210955608663024
$ws.Cells.Item(14, 17).Value = 0.5547002874985556
$ws.Cells.Item(14, 18).Value = 4.992302587487
$ws.Cells.Item(14, 19).Value = 0.004350743831597228
$ws.Cells.Item(14, 20).Value = 0.004350743831597228

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Tfpi"
$ws.Cells.Item(15, 3).Value = "Vldlr"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 1.799827666666667
$ws.Cells.Item(15, 8).Value = 5.399483
$ws.Cells.Item(15, 9).Value = 0.04723444576700921
$ws.Cells.Item(15, 10).Value = 0.0472344457670092
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 2.207039333333333
$ws.Cells.Item(15, 14).Value = 6.621118
$ws.Cells.Item(15, 15).Value = 0.65961009678592
$ws.Cells.Item(15, 16).Value = 0.6596100967859201
$ws.Cells.Item(15, 17).Value = 3.972290453554889
$ws.Cells.Item(15, 18).Value = 35.750614081994
$ws.Cells.Item(15, 19).Value = 0.03115631734400623
$ws.Cells.Item(15, 20).Value = 0.03115631734400623

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Tfpi"
$ws.Cells.Item(16, 3).Value = "Vldlr"
$ws.Cells.Item(16, 4).Value = "sCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 1.799827666666667
$ws.Cells.Item(16, 8).Value = 5.399483
$ws.Cells.Item(16, 9).Value = 0.04723444576700921
$ws.Cells.Item(16, 10).Value = 0.0472344457670092
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.83074
$ws.Cells.Item(16, 14).Value = 2.49222
$ws.Cells.Item(16, 15).Value = 0.2482803471274497
$ws.Cells.Item(16, 16).Value = 0.2482803471274497
$ws.Cells.Item(16, 17).Value = 1.495188835806667
$ws.Cells.Item(16, 18).Value = 13.45669952226
$ws.Cells.Item(16, 19).Value = 0.01172738459140574
$ws.Cells.Item(16, 20).Value = 0.01172738459140574
